$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.132.35"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.427.17"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "2.426.91"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "2.859.19"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "62.013.68"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "2.424.69"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("D28").Value = "2.544.70"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +0.24%  "
